# HR1 case study updated.
#
# This script reproduces the changes from the commit:
#  - Rename sheet "CostFlex" -> "CostFlex, Winter"
#  - Update the selection/active-cell state on several worksheets
#    ("Main", "Pg, Winter, S1", "GenStatus, Winter", "CostFlex, Winter")
#    and make "CostFlex, Winter" the final active/selected sheet
#  - Replace the simulated cost-flexibility data table (B2:Y32) on the
#    "CostFlex, Winter" sheet with the new result values

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename "CostFlex" to "CostFlex, Winter"
# ---------------------------------------------------------------------
$wsCost = $wb.Worksheets.Item("CostFlex")
$wsCost.Name = "CostFlex, Winter"

# ---------------------------------------------------------------------
# 2) Update selection on "Main": F6 -> A3:B7
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Activate()
$wsMain.Range("A3:B7").Select()

# ---------------------------------------------------------------------
# 3) Update selection on "Pg, Winter, S1": F6 (sqref A1:Y15) -> L31
# ---------------------------------------------------------------------
$wsPg = $wb.Worksheets.Item("Pg, Winter, S1")
$wsPg.Activate()
$wsPg.Range("L31").Select()

# ---------------------------------------------------------------------
# 4) Update selection on "GenStatus, Winter": B2 -> P16
#    (this also naturally drops the old tabSelected/topLeftCell state
#    since it is no longer the last-activated sheet)
# ---------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("GenStatus, Winter")
$wsGen.Activate()
$wsGen.Range("P16").Select()

# ---------------------------------------------------------------------
# 5) Replace the data table on "CostFlex, Winter" (B2:Y32).
#    Every row (2-32) gets the same 24 new column values.
# ---------------------------------------------------------------------
$colVals = @(
    14.186000000000002,
    15.770000000000001,
    11.728,
    13.672000000000001,
    13.756,
    13.479999999999999,
    15.692000000000002,
    13.834,
    9.8580000000000005,
    10.110000000000001,
    7.5760000000000005,
    7.24,
    7.5180000000000007,
    9.032,
    8.2260000000000009,
    8.7379999999999995,
    8.3099999999999987,
    7.7359999999999998,
    6.3180000000000005,
    3.95,
    4.7459999999999996,
    6.0860000000000003,
    8.0180000000000007,
    9.6480000000000015
)

$rowCount = 31   # rows 2 .. 32
$colCount = 24   # columns B .. Y

$data = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $data[$r, $c] = $colVals[$c]
    }
}

$wsCost.Range("B2:Y32").Value = $data

# ---------------------------------------------------------------------
# 6) Finally activate "CostFlex, Winter" and select I13, making it the
#    active sheet/tab (matches activeTab moving to this sheet).
# ---------------------------------------------------------------------
$wsCost.Activate()
$wsCost.Range("I13").Select()
